$wb = $excel.ActiveWorkbook

# --- 1. Insert new "Player Info" sheet as the first sheet ---
$firstSheet = $wb.Worksheets.Item(1)
$wsInfo = $wb.Worksheets.Add($firstSheet)
$wsInfo.Name = "Player Info"

$wsInfo.Cells.Item(1, 1).Value = "ID"
$wsInfo.Cells.Item(1, 2).Value = "NAME"
$wsInfo.Cells.Item(1, 3).Value = "BATTING_HAND"
$wsInfo.Cells.Item(1, 4).Value = "BOWL_STYLE"

# Re-use the existing bold/bordered header style from the "ODI Batting" sheet
$wb.Worksheets.Item("ODI Batting").Range("A1").Copy()
$wsInfo.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsInfo.Cells.Item(2, 1).Value = "'5925"
$wsInfo.Cells.Item(2, 2).Value = "Kyle Rico Mayers"
$wsInfo.Cells.Item(2, 3).Value = "Left Handed"
$wsInfo.Cells.Item(2, 4).Value = "Right Arm Medium"

$wsInfo.Range("A1").Select()

# --- 2. Update "ODI Batting" sheet: rename MATCH_CARD_LINK -> MATCH_CODE, store bare match codes ---
$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsBatting.Cells.Item(1, 4).Value = "MATCH_CODE"
$wsBatting.Cells.Item(2, 4).Value = "'4443"
$wsBatting.Cells.Item(3, 4).Value = "'4445"
$wsBatting.Cells.Item(4, 4).Value = "'4447"
$wsBatting.Cells.Item(5, 4).Value = "'4577"
$wsBatting.Cells.Item(6, 4).Value = "'4580"
$wsBatting.Cells.Item(7, 4).Value = "'4583"
$wsBatting.Cells.Item(8, 4).Value = "'4586"
$wsBatting.Cells.Item(9, 4).Value = "'4590"
$wsBatting.Cells.Item(10, 4).Value = "'4592"
$wsBatting.Cells.Item(11, 4).Value = "'4606"
$wsBatting.Cells.Item(12, 4).Value = "'4611"
$wsBatting.Cells.Item(13, 4).Value = "'4621"
$wsBatting.Cells.Item(14, 4).Value = "'4623"
$wsBatting.Cells.Item(15, 4).Value = "'4624"
$wsBatting.Cells.Item(16, 4).Value = "'4636"
$wsBatting.Cells.Item(17, 4).Value = "'4639"
$wsBatting.Cells.Item(18, 4).Value = "'4642"
$wsBatting.Cells.Item(19, 4).Value = "'4727"
$wsBatting.Cells.Item(20, 4).Value = "'4731"

# --- 3. Update "ODI Bowling" sheet: rename MATCH_CARD_LINK -> MATCH_CODE, store bare match codes ---
$wsBowling = $wb.Worksheets.Item("ODI Bowling")
$wsBowling.Cells.Item(1, 2).Value = "MATCH_CODE"
$wsBowling.Cells.Item(2, 2).Value = "'4445"
$wsBowling.Cells.Item(3, 2).Value = "'4447"
$wsBowling.Cells.Item(4, 2).Value = "'4577"
$wsBowling.Cells.Item(5, 2).Value = "'4580"
$wsBowling.Cells.Item(6, 2).Value = "'4583"
$wsBowling.Cells.Item(7, 2).Value = "'4586"
$wsBowling.Cells.Item(8, 2).Value = "'4590"
$wsBowling.Cells.Item(9, 2).Value = "'4621"
$wsBowling.Cells.Item(10, 2).Value = "'4623"
$wsBowling.Cells.Item(11, 2).Value = "'4624"
$wsBowling.Cells.Item(12, 2).Value = "'4639"
$wsBowling.Cells.Item(13, 2).Value = "'4642"
$wsBowling.Cells.Item(14, 2).Value = "'4727"
$wsBowling.Cells.Item(15, 2).Value = "'4731"

$wb.Worksheets.Item(1).Activate()
